$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index map: A=1 B=2 C=3 D=4 E=5

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.061.17"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -2.20%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.669.30"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -1.54%  "

$ws.Cells.Item(4, 5).Value = "  -0.23%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "216.81"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.42%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5107"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.38%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.21%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2656"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.46%  "

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06415"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.02%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.85"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.39%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07430"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.09%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.686.50"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.60%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.506"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.27%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5844"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.95%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.000008590"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.86%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "64.50"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.53%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "26.077.90"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.23%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.952"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.82%  "

$ws.Cells.Item(19, 5).Value = "  -0.12%  "

$ws.Cells.Item(20, 5).Value = "  -1.72%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "193.29"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.52%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.228"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.36%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.24%  "

$ws.Cells.Item(24, 5).Value = "  +0.22%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.617"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.53%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1196"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +3.28%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.71"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.64%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06413"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +13.38%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.340"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.30%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.319"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.28%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.555"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.36%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.523"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.95%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.650"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.09%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.020"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.03%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6112"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.80%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.368"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.26%  "

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.706"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.36%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.259"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +7.12%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01604"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.82%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.089.84"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.22%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.8626"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.38%  "

$ws.Cells.Item(42, 5).Value = "  +0.55%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "100.75"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.14%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.817.42"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.93%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00000000113"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.51%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "56.47"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.11%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.23%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.074"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.24%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05237"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.05%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4285"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.90%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.036"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +4.28%  "
